# Applies the "Anonimyzed fedcore" update:
#  - In both sheets, give the two cells to the right of the first merged
#    header cell (e.g. C1/D1, and F1/G1) a thin top+bottom box border
#    (matching the already-present border definitions: index 4 is
#    "top+bottom" and index 5 is "top+right+bottom", used elsewhere in
#    the workbook's styles).
#  - Rename every occurrence of "fedcore" in row 2 headers to "approach".
#  - Remove the stray empty cell G5 on the computational_comparison sheet.

$wb = $excel.ActiveWorkbook

function Set-HeaderBoxBorder($cell, [bool]$withRightEdge) {
    # Start from a clean slate so the resulting style only carries the
    # default font/fill plus the border we are about to add (this mirrors
    # fontId="0" / no alignment in the target style records).
    $cell.ClearFormats()

    $cell.Borders.Item(8).LineStyle = 1       # xlEdgeTop    -> continuous
    $cell.Borders.Item(8).Weight = 2          # xlThin
    $cell.Borders.Item(9).LineStyle = 1       # xlEdgeBottom -> continuous
    $cell.Borders.Item(9).Weight = 2
    if ($withRightEdge) {
        $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight  -> continuous
        $cell.Borders.Item(10).Weight = 2
    }
}

# ---- Sheet "quality_comparison" ----
$ws1 = $wb.Worksheets.Item("quality_comparison")

$ws1C1 = $ws1.Range("C1")
$ws1D1 = $ws1.Range("D1")
Set-HeaderBoxBorder $ws1C1 $false
Set-HeaderBoxBorder $ws1D1 $true

$ws1.Range("C2").Value = "approach"

# ---- Sheet "computational_comparison" ----
$ws2 = $wb.Worksheets.Item("computational_comparison")

$ws2C1 = $ws2.Range("C1")
$ws2D1 = $ws2.Range("D1")
Set-HeaderBoxBorder $ws2C1 $false
Set-HeaderBoxBorder $ws2D1 $true

# F1/G1 need the exact same resulting styles as C1/D1. Setting the border
# edges one property at a time a second time would allocate fresh (and
# then orphaned) style records, so instead copy the already-resolved
# formatting over with copy/paste-special, which reuses the existing
# style entries.
$ws2C1.Copy()
$ws2.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws2D1.Copy()
$ws2.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5 entirely.
$ws2.Range("G5").ClearContents()
